$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.286965847015381
$ws.Range("B1").Value = 3.94300365447998
$ws.Range("C1").Value = 3.691750049591064
$ws.Range("D1").Value = 3.278568983078003
$ws.Range("E1").Value = 1.053878784179688
